# Updated cryptos list on Fri Mar 24 20:46:04 UTC 2023 with GitHub Actions
#
# This script applies the latest scraped price/volume snapshot to the
# cryptocurrency table on the active worksheet. Column D ("Price") holds
# locale-formatted numeric-looking strings (e.g. "1.780.87", "6.040") that
# must remain plain text -- they are prefixed with a leading apostrophe so
# Excel stores them as text instead of re-parsing/rounding them as numbers.
# Column E ("Volume(1h)") values already carry padding spaces and a percent
# sign, so Excel keeps them as text natively.
# Two coin pairs (rows 39/40 and 46/47) swapped rank order in this refresh,
# so their Coin name (B) and Link (C) cells are updated along with price/volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.702.65"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "'1.753.56"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").Value = "'323.78"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("E7").Value = "  -3.96%  "
$ws.Range("D8").Value = "'0.3623"
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("D9").Value = "'0.07563"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").Value = "'42.51"
$ws.Range("E10").Value = "  -5.26%  "
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "'20.62"
$ws.Range("E13").Value = "  -6.77%  "
$ws.Range("D14").Value = "'6.040"
$ws.Range("E14").Value = "  -3.50%  "
$ws.Range("D15").Value = "'7.261"
$ws.Range("E15").Value = "  -4.17%  "
$ws.Range("D16").Value = "'1.780.87"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "'91.16"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "'0.00001074"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "'0.06377"
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "'17.01"
$ws.Range("E21").Value = "  -2.99%  "
$ws.Range("D22").Value = "'5.899"
$ws.Range("E22").Value = "  -5.08%  "
$ws.Range("D23").Value = "'27.778.23"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("D24").Value = "'11.19"
$ws.Range("E24").Value = "  -4.39%  "
$ws.Range("D25").Value = "'2.099"
$ws.Range("E25").Value = "  +7.21%  "
$ws.Range("D26").Value = "'160.17"
$ws.Range("E26").Value = "  +2.85%  "
$ws.Range("D27").Value = "'20.29"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").Value = "'1.974.01"
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").Value = "'2.133"
$ws.Range("E29").Value = "  -8.17%  "
$ws.Range("D30").Value = "'125.18"
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("D31").Value = "'1.114"
$ws.Range("E31").Value = "  -7.44%  "
$ws.Range("D32").Value = "'3.679"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "'5.570"
$ws.Range("E33").Value = "  -5.34%  "
$ws.Range("D34").Value = "'0.08887"
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("D35").Value = "'12.24"
$ws.Range("E35").Value = "  -6.77%  "
$ws.Range("D36").Value = "'0.02299"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").Value = "'0.2103"
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("D38").Value = "'0.06015"
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.6337"
$ws.Range("E39").Value = "  -3.72%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'4.961"
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").Value = "'7.898"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("D44").Value = "'1.393"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "'13.28"
$ws.Range("E45").Value = "  -4.51%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.698"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5861"
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("D48").Value = "'1.985"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("D49").Value = "'122.89"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").Value = "'1.178"
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("D51").Value = "'0.06821"
$ws.Range("E51").Value = "  -2.39%  "
